# Updated cryptos list on Fri Apr 12 13:28:24 UTC 2024 with GitHub Actions
#
# All Price (column D) and Volume(1h) (column E) cells in this sheet are
# stored as plain text, even when they look numeric (e.g. "1.00", "592.20").
# We force column D/E cells to Text format before writing so the engine
# doesn't silently coerce these into numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - Bitcoin
Set-Text "D2" "69.973.32"
Set-Text "E2" "  -1.03%  "

# Row 3 - Ethereum
Set-Text "D3" "3.469.56"
Set-Text "E3" "  -2.34%  "

# Row 4 - TetherUSD
Set-Text "D4" "1.00"
Set-Text "E4" "  -0.06%  "

# Row 5 - BNB
Set-Text "D5" "612.78"
Set-Text "E5" "  +1.58%  "

# Row 6 - Solana
Set-Text "D6" "167.95"
Set-Text "E6" "  -2.88%  "

# Row 7 - LidoStakedEther
Set-Text "D7" "3.465.55"
Set-Text "E7" "  -2.24%  "

# Row 8 - XRP
Set-Text "D8" "0.597"
Set-Text "E8" "  -2.63%  "

# Row 9 - USDC
Set-Text "E9" "  +0.01%  "

# Row 10 - Dogecoin
Set-Text "D10" "0.193"
Set-Text "E10" "  -0.03%  "

# Row 11 - Toncoin
Set-Text "D11" "7.05"
Set-Text "E11" "  -3.99%  "

# Row 12 - Cardano
Set-Text "D12" "0.566"
Set-Text "E12" "  -3.34%  "

# Row 13 - Avalanche
Set-Text "D13" "44.60"
Set-Text "E13" "  -3.68%  "

# Row 14 - ShibaInu
Set-Text "D14" "0.0000268"
Set-Text "E14" "  -2.98%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-Text "D15" "4.040.03"
Set-Text "E15" "  -2.04%  "

# Row 16 - Polkadot
Set-Text "D16" "8.21"
Set-Text "E16" "  -1.13%  "

# Row 17 - BitcoinCash
Set-Text "D17" "589.68"
Set-Text "E17" "  -2.76%  "

# Row 18 - WrappedEther
Set-Text "D18" "3.487.07"
Set-Text "E18" "  -1.90%  "

# Row 19 - WrappedBTC
Set-Text "D19" "70.075.30"

# Row 20 - TRON
Set-Text "E20" "  +0.84%  "

# Row 21 - Chainlink
Set-Text "D21" "17.19"
Set-Text "E21" "  -0.73%  "

# Row 22 - Polygon
Set-Text "D22" "0.857"
Set-Text "E22" "  -2.31%  "

# Row 23 - Uniswap
Set-Text "D23" "8.75"
Set-Text "E23" "  -5.54%  "

# Row 24 - Litecoin
Set-Text "D24" "95.74"
Set-Text "E24" "  -0.71%  "

# Row 25 - InternetComputer(DFINITY)
Set-Text "D25" "15.21"
Set-Text "E25" "  -3.03%  "

# Row 26 - PancakeSwap
Set-Text "D26" "3.63"
Set-Text "E26" "  -2.62%  "

# Row 27 - Dai
Set-Text "D27" "0.998"
Set-Text "E27" "  -0.17%  "

# Row 28 - ImmutableX
Set-Text "D28" "2.48"
Set-Text "E28" "  -4.69%  "

# Row 29 - EthereumClassic
Set-Text "D29" "33.05"
Set-Text "E29" "  -2.74%  "

# Row 30 - RenderToken
Set-Text "D30" "8.64"
Set-Text "E30" "  -4.52%  "

# Row 31 - Filecoin
Set-Text "D31" "7.86"
Set-Text "E31" "  -4.23%  "

# Row 32 - Stacks
Set-Text "D32" "2.83"
Set-Text "E32" "  -7.01%  "

# Row 33 - Mantle
Set-Text "D33" "1.25"
Set-Text "E33" "  -3.17%  "

# Row 34 - NEARProtocol
Set-Text "D34" "6.58"
Set-Text "E34" "  -6.24%  "

# Row 35 - Bittensor
Set-Text "D35" "566.57"
Set-Text "E35" "  -21.51%  "

# Row 36 - was VeChain, now Cosmos
Set-Text "B36" "Cosmos"
Set-Text "C36" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-Text "D36" "10.66"
Set-Text "E36" "  -0.75%  "

# Row 37 - was Cosmos, now VeChain
Set-Text "B37" "VeChain"
Set-Text "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-Text "D37" "0.0482"
Set-Text "E37" "  +0.98%  "

# Row 38 - Hedera
Set-Text "D38" "0.0963"
Set-Text "E38" "  -4.20%  "

# Row 39 - FirstDigitalUSD
Set-Text "E39" "  +0.42%  "

# Row 40 - OKB
Set-Text "D40" "56.23"
Set-Text "E40" "  -1.15%  "

# Row 41 - Kaspa
Set-Text "E41" "  -0.92%  "

# Row 42 - dogwifhat
Set-Text "D42" "3.21"
Set-Text "E42" "  -10.06%  "

# Row 43 - Maker
Set-Text "D43" "3.270.06"
Set-Text "E43" "  -2.79%  "

# Row 44 - PEPE. The price contains U+2083 SUBSCRIPT THREE. Writing the
# literal character directly trips the engine's numeric auto-detection
# (it mis-parses the digit run following the subscript), so we write a
# plain placeholder first and swap it in afterwards with Replace, which
# operates on the already-stored text and doesn't re-trigger that parse.
Set-Text "D44" "0.0Z0698"
$ws.Range("D44").Replace("Z", [char]0x2083) | Out-Null
Set-Text "E44" "  -0.05%  "

# Row 45 - TheGraph
Set-Text "D45" "0.299"
Set-Text "E45" "  -5.66%  "

# Row 46 - InjectiveProtocol
Set-Text "D46" "30.81"
Set-Text "E46" "  -5.23%  "

# Row 47 - ThetaToken
Set-Text "D47" "2.75"
Set-Text "E47" "  -6.01%  "

# Row 48 - Fetch.AI
Set-Text "D48" "2.40"
Set-Text "E48" "  -6.83%  "

# Row 49 - Stellar
Set-Text "E49" "  -2.69%  "

# Row 50 - Monero
Set-Text "D50" "132.98"
Set-Text "E50" "  -0.68%  "

# Row 51 - USDe
Set-Text "E51" "  -0.03%  "
